# Author's edit: add the "NAAN MUDHALVAN ID" line to the details textbox
# on the title slide (slide 1), right after the "NAAN MUDHALVAN
# USERNAME:..." line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the details textbox by name ("TextBox 13") so this is robust to
# shape ordering.
$detailsShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 13") {
        $detailsShape = $candidate
    }
}

$tr = $detailsShape.TextFrame.TextRange

# Paragraph 3 is "NAAN MUDHALVAN USERNAME:asunm110312201163"; insert a new
# paragraph straight after it containing the NAAN MUDHALVAN ID line. The
# shape auto-fits its height (spAutoFit), so it grows to accommodate the
# extra line automatically.
$usernamePara = $tr.Paragraphs(3, 1)
$null = $usernamePara.InsertAfter("`rNAAN MUDHALVAN ID:301A25241C15348A7ED4E95FF6A2D40E")
